$d = $word.ActiveDocument

function Set-ParaXML($para, $innerBodyXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $innerBodyXml +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $para.Range.InsertXML($xml)
}

$spacingPPr = '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'

# --- 1. "To Lunaura, land of cyan skies," (Verse 1) ---
$target = $d.Paragraphs.Item(2)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">To </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Lunaura</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, land of cyan skies,</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner

# --- 2. First Chorus: "Oh Lunaura, sacred home," / "Oh Lunaura, where we've grown," ---
$target = $d.Paragraphs.Item(8)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">Oh </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Lunaura</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, sacred home,</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner

$target = $d.Paragraphs.Item(9)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">Oh </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Lunaura</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, where we’ve grown,</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner

# --- 3. Second Chorus (after Verse 2): same two lines ---
$target = $d.Paragraphs.Item(20)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">Oh </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Lunaura</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, sacred home,</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner

$target = $d.Paragraphs.Item(21)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">Oh </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Lunaura</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, where we’ve grown,</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner

# --- 4. Delete the "(Verse 3)" section and the Chorus repeat that follows it ---
# At this point paragraph numbering is unchanged (edits above kept 1 paragraph -> 1 paragraph).
$startPara = $d.Paragraphs.Item(25)
$endPara = $d.Paragraphs.Item(36)
$rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)
$null = $rangeToDelete.Delete()

# --- 5. "(Outro)" block: "For Lunaura, land divine," is now paragraph 26 ---
$target = $d.Paragraphs.Item(26)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">For </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Lunaura</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, land divine,</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner

# --- 6. "We stand by you, til the last of days." is now paragraph 29 ---
$target = $d.Paragraphs.Item(29)
$inner = '<w:body><w:p>' + $spacingPPr +
    '<w:r><w:t xml:space="preserve">We stand by you, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>til</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> the last of days</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p></w:body>'
Set-ParaXML $target $inner
